$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Apply per-cell updates from the weekly refresh (reshuffled rows 2-32) ---
$ws.Range("D3").Value = 44266
$ws.Range("J3").Value = 160
$ws.Range("D4").Value = 44279
$ws.Range("K4").Value = 10000
$ws.Range("M4").Value = 10000
$ws.Range("P4").Value = 3333
$ws.Range("D6").Value = 44342
$ws.Range("J6").Value = 17
$ws.Range("L6").Value = 10000
$ws.Range("M6").Value = 10000
$ws.Range("P6").Value = 3333
$ws.Range("D7").Value = 44405
$ws.Range("L7").Value = 10500
$ws.Range("M7").Value = 10250
$ws.Range("P7").Value = 3417
$ws.Range("D8").Value = 44272
$ws.Range("J8").Value = 70
$ws.Range("D9").Value = 44419
$ws.Range("J9").Value = 16
$ws.Range("L9").Value = 10000
$ws.Range("M9").Value = 10000
$ws.Range("P9").Value = 3333
$ws.Range("D10").Value = 44370
$ws.Range("D11").Value = 44412
$ws.Range("J11").Value = 25
$ws.Range("M11").Value = 10260
$ws.Range("P11").Value = 3420
$ws.Range("D12").Value = 44293
$ws.Range("J12").Value = 16
$ws.Range("D13").Value = 44195
$ws.Range("J13").Value = 30
$ws.Range("D14").Value = 44447
$ws.Range("L14").Value = 10500
$ws.Range("M14").Value = 10250
$ws.Range("P14").Value = 3417
$ws.Range("D15").Value = 44356
$ws.Range("D16").Value = 44454
$ws.Range("K16").Value = 9500
$ws.Range("M16").Value = 9750
$ws.Range("P16").Value = 3250
$ws.Range("D17").Value = 44426
$ws.Range("J17").Value = 16
$ws.Range("L17").Value = 10500
$ws.Range("M17").Value = 10250
$ws.Range("P17").Value = 3417
$ws.Range("D18").Value = 44391
$ws.Range("J18").Value = 16
$ws.Range("D19").Value = 44335
$ws.Range("D20").Value = 44475
$ws.Range("K20").Value = 9000
$ws.Range("L20").Value = 10000
$ws.Range("M20").Value = 9500
$ws.Range("P20").Value = 3167
$ws.Range("D21").Value = 44300
$ws.Range("J21").Value = 16
$ws.Range("D22").Value = 44349
$ws.Range("J22").Value = 12
$ws.Range("L22").Value = 10000
$ws.Range("M22").Value = 10000
$ws.Range("P22").Value = 3333
$ws.Range("D23").Value = 44435
$ws.Range("L23").Value = 10500
$ws.Range("M23").Value = 10250
$ws.Range("P23").Value = 3417
$ws.Range("D24").Value = 44377
$ws.Range("D25").Value = 44433
$ws.Range("J25").Value = 16
$ws.Range("L25").Value = 10500
$ws.Range("M25").Value = 10250
$ws.Range("P25").Value = 3417
$ws.Range("D26").Value = 44307
$ws.Range("J26").Value = 160
$ws.Range("L26").Value = 10000
$ws.Range("M26").Value = 10000
$ws.Range("P26").Value = 3333
$ws.Range("D28").Value = 44384
$ws.Range("J28").Value = 25
$ws.Range("M28").Value = 10260
$ws.Range("P28").Value = 3420
$ws.Range("D29").Value = 44363
$ws.Range("D30").Value = 44328
$ws.Range("D31").Value = 44321
$ws.Range("J31").Value = 25
$ws.Range("D32").Value = 44314
$ws.Range("K32").Value = 10000
$ws.Range("L32").Value = 10000
$ws.Range("M32").Value = 10000
$ws.Range("P32").Value = 3333

# --- Append new row 33 (new weekly data point; old row 32 tail shifted down) ---
$ws.Range("A33").Value = 9
$ws.Range("B33").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C33").Value = "Metropolitana"
$ws.Range("D33").Value = 44461
$ws.Range("E33").Value = 13
$ws.Range("F33").Value = 100112029
$ws.Range("G33").Value = "Orégano"
$ws.Range("H33").Value = "Sin especificar"
$ws.Range("I33").Value = "Primera"
$ws.Range("J33").Value = 16
$ws.Range("K33").Value = 9500
$ws.Range("L33").Value = 10000
$ws.Range("M33").Value = 9750
$ws.Range("N33").Value = "$/docena de atados"
$ws.Range("O33").Value = "Región Metropolitana"
$ws.Range("P33").Value = 3250
$ws.Range("Q33").Value = 3
$ws.Range("R33").Value = "Hortaliza"

# Match the date style/number-format used by the D column on existing rows
$ws.Range("D33").NumberFormat = $ws.Range("D32").NumberFormat
